# Auto-generated Excel COM-interop script that applies the
# Marilith_Profits market-data refresh described in the commit diff.
# For every (sheet, row, column) touched by the diff we either:
#   - write the new numeric value (modify / add), or
#   - clear the cell entirely (remove), matching the OOXML diff
#     semantics where the <c> element disappears rather than
#     becoming 0 / blank.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 194.5
$ws.Range("I9").Value = 190
$ws.Range("J9").Value = 199
$ws.Range("K9").Value = 190
$ws.Range("L9").Value = 199
$ws.Range("M9").Value = -21
$ws.Range("N9").Value = -537
$ws.Range("H15").Value = 1348.9565
$ws.Range("I15").Value = 1348.9565
$ws.Range("K15").Value = 4046.8695
$ws.Range("M15").Value = -3877.8695
$ws.Range("H28").Value = 1105.1111
$ws.Range("I28").Value = 416.16666
$ws.Range("K28").Value = 416.16666
$ws.Range("M28").Value = 68.83334000000002
$ws.Range("H74").Value = 59110.89
$ws.Range("I74").Value = 3999.75
$ws.Range("J74").Value = 500000
$ws.Range("K74").Value = 3999.75
$ws.Range("L74").Value = 500000
$ws.Range("M74").Value = -3063.75
$ws.Range("N74").Value = -501872
$ws.Range("H77").Value = 59110.89
$ws.Range("I77").Value = 3999.75
$ws.Range("J77").Value = 500000
$ws.Range("K77").Value = 19998.75
$ws.Range("L77").Value = 2500000
$ws.Range("M77").Value = -15318.75
$ws.Range("N77").Value = -2509360
$ws.Range("H107").Value = 782.55554
$ws.Range("I107").Value = 782.55554
$ws.Range("K107").Value = 782.55554
$ws.Range("M107").Value = 1137.44446
$ws.Range("H116").Value = 6499.3335
$ws.Range("I116").Value = 5000
$ws.Range("K116").Value = 5000
$ws.Range("M116").Value = -1558
$ws.Range("H132").Value = 1645.5358
$ws.Range("I132").Value = 1645.5358
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4936.607400000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2406.607400000001
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 725
$ws.Range("H137").Value = 3366.4375
$ws.Range("I137").Value = 1774.6
$ws.Range("K137").Value = 5323.799999999999
$ws.Range("M137").Value = -2773.799999999999
$ws.Range("H138").Value = 2928.5715
$ws.Range("J138").Value = 3333.3333
$ws.Range("L138").Value = 9999.999899999999
$ws.Range("N138").Value = -20279.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5389.9
$ws.Range("I74").Value = 4872.5
$ws.Range("K74").Value = 4872.5
$ws.Range("M74").Value = -3998.5
$ws.Range("H77").Value = 5389.9
$ws.Range("I77").Value = 4872.5
$ws.Range("K77").Value = 24362.5
$ws.Range("M77").Value = -19994.5
$ws.Range("H122").Value = 997.8333
$ws.Range("I122").Value = 997.8333
$ws.Range("K122").Value = 2993.4999
$ws.Range("M122").Value = -543.4998999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6996
$ws.Range("I20").Value = 6997.3335
$ws.Range("J20").Value = 6994
$ws.Range("K20").Value = 6997.3335
$ws.Range("L20").Value = 6994
$ws.Range("M20").Value = -6750.3335
$ws.Range("N20").Value = -7488
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H94").Value = 2856
$ws.Range("I94").Value = 2856
$ws.Range("K94").Value = 2856
$ws.Range("M94").Value = -2405

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 44.454544
$ws.Range("J7").Value = 39.25
$ws.Range("L7").Value = 39.25
$ws.Range("N7").Value = -265.25
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -650
$ws.Range("N22").ClearContents()
$ws.Range("H107").Value = 757
$ws.Range("I107").Value = 709.8
$ws.Range("K107").Value = 709.8
$ws.Range("M107").Value = 1210.2
$ws.Range("H132").Value = 870.4091
$ws.Range("I132").Value = 870.4091
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2611.2273
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -81.22730000000001
$ws.Range("N132").ClearContents()
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 1771.3667
$ws.Range("I134").Value = 1771.3667
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5314.1001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2779.1001
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 825
$ws.Range("J16").Value = 970
$ws.Range("L16").Value = 2910
$ws.Range("N16").Value = -3256
$ws.Range("H29").Value = 35.25
$ws.Range("J29").Value = 23.809525
$ws.Range("L29").Value = 71.428575
$ws.Range("N29").Value = -625.428575
$ws.Range("H34").Value = 509.4
$ws.Range("I34").Value = 382.33334
$ws.Range("J34").Value = 700
$ws.Range("K34").Value = 1147.00002
$ws.Range("L34").Value = 2100
$ws.Range("M34").Value = -1063.00002
$ws.Range("N34").Value = -2268
$ws.Range("H38").Value = 235.66667
$ws.Range("I38").Value = 233.6
$ws.Range("K38").Value = 700.8
$ws.Range("M38").Value = -353.8
$ws.Range("H39").Value = 4857.143
$ws.Range("J39").Value = 4857.143
$ws.Range("L39").Value = 14571.429
$ws.Range("N39").Value = -15159.429
$ws.Range("H46").Value = 1112454.8
$ws.Range("J46").Value = 1717.7142
$ws.Range("L46").Value = 5153.142599999999
$ws.Range("N46").Value = -5335.142599999999
$ws.Range("H55").Value = 3222.2222
$ws.Range("J55").Value = 5000
$ws.Range("L55").Value = 15000
$ws.Range("N55").Value = -15354
$ws.Range("H59").Value = 605
$ws.Range("I59").Value = 605
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 1815
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -1275
$ws.Range("N59").ClearContents()
$ws.Range("H107").Value = 989.8
$ws.Range("J107").Value = 855.875
$ws.Range("L107").Value = 2567.625
$ws.Range("N107").Value = -6407.625
$ws.Range("H122").Value = 734.75
$ws.Range("I122").Value = 1116.3334
$ws.Range("J122").Value = 607.55554
$ws.Range("K122").Value = 10047.0006
$ws.Range("L122").Value = 5467.99986
$ws.Range("M122").Value = -7597.000599999999
$ws.Range("N122").Value = -10367.99986
$ws.Range("H137").Value = 1130
$ws.Range("I137").Value = 1130
$ws.Range("K137").Value = 3390
$ws.Range("M137").Value = 1710

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 508.33334
$ws.Range("I102").Value = 521.875
$ws.Range("K102").Value = 521.875
$ws.Range("M102").Value = 1100.125
$ws.Range("H122").Value = 3541.3635
$ws.Range("I122").Value = 2349.8572
$ws.Range("J122").Value = 5626.5
$ws.Range("K122").Value = 7049.571599999999
$ws.Range("L122").Value = 16879.5
$ws.Range("M122").Value = -4599.571599999999
$ws.Range("N122").Value = -21779.5
$ws.Range("H132").Value = 3995.3333
$ws.Range("I132").Value = 3995.3333
$ws.Range("K132").Value = 11985.9999
$ws.Range("M132").Value = -9455.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3049.6
$ws.Range("I122").Value = 2937
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 8811
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -6361
$ws.Range("N122").Value = -15400
$ws.Range("H136").Value = 3161.2307
$ws.Range("I136").Value = 3091.75
$ws.Range("J136").Value = 3995
$ws.Range("K136").Value = 9275.25
$ws.Range("L136").Value = 11985
$ws.Range("M136").Value = -6725.25
$ws.Range("N136").Value = -17085

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 3633.3333
$ws.Range("I28").Value = 2000
$ws.Range("J28").Value = 4450
$ws.Range("K28").Value = 2000
$ws.Range("L28").Value = 4450
$ws.Range("M28").Value = -1652
$ws.Range("N28").Value = -5146
$ws.Range("H33").Value = 12673.667
$ws.Range("J33").Value = 9010.5
$ws.Range("L33").Value = 9010.5
$ws.Range("N33").Value = -9510.5
$ws.Range("H36").Value = 12673.667
$ws.Range("J36").Value = 9010.5
$ws.Range("L36").Value = 9010.5
$ws.Range("N36").Value = -9510.5
$ws.Range("H54").Value = 21000
$ws.Range("I54").Value = 12000
$ws.Range("K54").Value = 12000
$ws.Range("M54").Value = -11480
$ws.Range("H122").Value = 1525.5
$ws.Range("I122").Value = 1412.0769
$ws.Range("K122").Value = 4236.2307
$ws.Range("M122").Value = -1786.2307
$ws.Range("H126").Value = 1182.7778
$ws.Range("I126").Value = 1142.2858
$ws.Range("K126").Value = 3426.8574
$ws.Range("M126").Value = -956.8574000000003
$ws.Range("H132").Value = 1939.5333
$ws.Range("I132").Value = 1939.5333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5818.5999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3288.5999
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 2102.6428
$ws.Range("I136").Value = 1576.5454
$ws.Range("K136").Value = 4729.6362
$ws.Range("M136").Value = -2179.6362

